$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - column headers
$ws.Range("A1").Value = "Municipio estudio, nombre"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Comarca nombre"
$ws.Range("D1").Value = "Municipio estudio, código"
$ws.Range("E1").Value = "Comarca código"
$ws.Range("F1").Value = "Provincia residencia código"
$ws.Range("G1").Value = "Municipio residencia nombre"
$ws.Range("H1").Value = "Municipio residencia código"
$ws.Range("I1").Value = "Provincia residencia nombre"

# Row 2 - measure/dimension identifiers
$ws.Range("A2").Value = "iaest-measure:municipio-estudio-nombre"
$ws.Range("B2").Value = "iaest-measure:personas-residentes-viviendas-familiares"
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "null"
$ws.Range("F2").Value = "null"
$ws.Range("G2").Value = "iaest-measure:municipio-residencia-nombre"
$ws.Range("H2").Value = "null"
$ws.Range("I2").Value = "iaest-measure:provincia-residencia-nombre"

# Row 3 - medida/dim markers
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "dim"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "null"
$ws.Range("I3").Value = "medida"

# Row 4 - xsd types / codelist reference
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("C4").Value = "URI-comarca"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "null"
$ws.Range("G4").Value = "xsd:string"
$ws.Range("H4").Value = "null"
$ws.Range("I4").Value = "xsd:string"
